# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
#
# The upstream source re-pulled odds data and a handful of rows ended up
# re-ordered (their B:AC payload moved to a different row while the leading
# row-index column A stayed put) plus a few rows got refreshed odds values
# in place.
#
# Columns used on this sheet: A..AC (1..29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colFirst = 2   # column B
$colLast  = 29  # column AC

function Swap-RowData($rowA, $rowB) {
    for ($c = $colFirst; $c -le $colLast; $c++) {
        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

function Rotate-RowData($row1, $row2, $row3) {
    # new(row1) = old(row2); new(row2) = old(row3); new(row3) = old(row1)
    for ($c = $colFirst; $c -le $colLast; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $cell3 = $ws.Cells.Item($row3, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $v3 = $cell3.Value2
        $cell1.Value = $v2
        $cell2.Value = $v3
        $cell3.Value = $v1
    }
}

# Rows 95 / 96 swapped (same match ids, re-ordered).
Swap-RowData 95 96

# Rows 129 / 130 / 131 rotated: 129<-130, 130<-131, 131<-129.
Rotate-RowData 129 130 131

# Rows 224 / 225 swapped.
Swap-RowData 224 225

# Row 230: refreshed odds values (no row reorder).
$ws.Cells.Item(230, 14).Value = 1.666    # N230
$ws.Cells.Item(230, 15).Value = 3.6      # O230
$ws.Cells.Item(230, 16).Value = 5.25     # P230
$ws.Cells.Item(230, 17).Value = -0.75    # Q230
$ws.Cells.Item(230, 18).Value = 1.9      # R230
$ws.Cells.Item(230, 19).Value = 1.9      # S230
$ws.Cells.Item(230, 21).Value = 1.85     # U230
$ws.Cells.Item(230, 22).Value = 1.95     # V230

# Row 233: refreshed odds values (no row reorder).
$ws.Cells.Item(233, 14).Value = 2.55     # N233
$ws.Cells.Item(233, 16).Value = 2.7      # P233
$ws.Cells.Item(233, 18).Value = 1.825    # R233
$ws.Cells.Item(233, 19).Value = 1.975    # S233

# Row 235: refreshed odds values (no row reorder).
$ws.Cells.Item(235, 18).Value = 2.025    # R235
$ws.Cells.Item(235, 19).Value = 1.775    # S235
$ws.Cells.Item(235, 21).Value = 1.95     # U235
$ws.Cells.Item(235, 22).Value = 1.85     # V235

Write-Output "Applied league base update."
